$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 8.551
$ws.Range("D4").Value = -7.976000000000001
$ws.Range("D7").Value = -7.813999999999998
$ws.Range("A9").Value = -21.723
$ws.Range("B9").Value = 5.554
$ws.Range("C9").Value = -10.914
$ws.Range("D11").Value = -7.525
$ws.Range("D15").Value = -8.098000000000001
$ws.Range("A18").Value = -21.649
$ws.Range("A20").Value = -20.265
$ws.Range("B23").Value = 7.463000000000001
$ws.Range("B24").Value = 5.394
$ws.Range("B26").Value = 5.855
$ws.Range("A27").Value = -21.875
$ws.Range("D30").Value = -7.439
$ws.Range("C32").Value = -12.324
$ws.Range("B34").Value = 7.279999999999999
$ws.Range("B35").Value = 8.031000000000001
$ws.Range("C38").Value = -12.537
$ws.Range("D39").Value = -7.553
$ws.Range("D43").Value = -7.312
$ws.Range("C45").Value = -13.478
$ws.Range("D47").Value = -7.203999999999999
$ws.Range("B48").Value = 5.456999999999999
$ws.Range("C51").Value = -11.584
$ws.Range("B52").Value = 5.01
$ws.Range("C57").Value = -13.672
$ws.Range("C64").Value = -10.754
$ws.Range("B66").Value = 5.063
$ws.Range("B67").Value = 5.392
$ws.Range("A69").Value = -21.507
$ws.Range("D75").Value = -8.021000000000001
$ws.Range("A76").Value = -20.574
$ws.Range("B80").Value = 7.694
$ws.Range("A82").Value = -21.749
$ws.Range("D91").Value = -7.031999999999999
$ws.Range("D92").Value = -7.036
$ws.Range("C93").Value = -10.75
$ws.Range("B99").Value = 4.968999999999999
